# Update NATMI LR-pair TPM values (Ccl2-Ccr5) per updated computation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.719288
$ws.Range("H2").Value = 11.157864
$ws.Range("I2").Value = 0.04235839908674209
$ws.Range("J2").Value = 0.04235839908674209
$ws.Range("M2").Value = 0.007957
$ws.Range("N2").Value = 0.023871
$ws.Range("O2").Value = 0.0002448939493579708
$ws.Range("P2").Value = 0.0002448939493579708
$ws.Range("Q2").Value = 0.029594374616
$ws.Range("R2").Value = 0.266349371544
$ws.Range("S2").Value = 0.00001037331564083333
$ws.Range("T2").Value = 0.00001037331564083333
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.719288
$ws.Range("H3").Value = 11.157864
$ws.Range("I3").Value = 0.04235839908674209
$ws.Range("J3").Value = 0.04235839908674209
$ws.Range("O3").Value = 0.003249135679578298
$ws.Range("P3").Value = 0.003249135679578299
$ws.Range("Q3").Value = 0.3926439943973334
$ws.Range("R3").Value = 3.533795949576
$ws.Range("S3").Value = 0.0001376281858025505
$ws.Range("T3").Value = 0.0001376281858025506
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.719288
$ws.Range("H4").Value = 11.157864
$ws.Range("I4").Value = 0.04235839908674209
$ws.Range("J4").Value = 0.04235839908674209
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02449766666666667
$ws.Range("N4").Value = 0.073493
$ws.Range("O4").Value = 0.0007539688752111494
$ws.Range("P4").Value = 0.0007539688752111494
$ws.Range("Q4").Value = 0.09111387766133335
$ws.Range("R4").Value = 0.820024898952
$ws.Range("S4").Value = 0.00003193691451517591
$ws.Range("T4").Value = 0.00003193691451517591
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.719288
$ws.Range("H5").Value = 11.157864
$ws.Range("I5").Value = 0.04235839908674209
$ws.Range("J5").Value = 0.04235839908674209
$ws.Range("M5").Value = 32.353591
$ws.Range("N5").Value = 97.060773
$ws.Range("O5").Value = 0.9957520014958525
$ws.Range("P5").Value = 0.9957520014958525
$ws.Range("Q5").Value = 120.332322763208
$ws.Range("R5").Value = 1082.990904868872
$ws.Range("S5").Value = 0.04217846067078353
$ws.Range("T5").Value = 0.04217846067078353
# Row 6
$ws.Range("I6").Value = 0.2979256989470644
$ws.Range("J6").Value = 0.2979256989470644
$ws.Range("M6").Value = 0.007957
$ws.Range("N6").Value = 0.023871
$ws.Range("O6").Value = 0.0002448939493579708
$ws.Range("P6").Value = 0.0002448939493579708
$ws.Range("Q6").Value = 0.2081505659436667
$ws.Range("R6").Value = 1.873355093493
$ws.Range("S6").Value = 0.00007296020103038043
$ws.Range("T6").Value = 0.00007296020103038043
# Row 7
$ws.Range("I7").Value = 0.2979256989470644
$ws.Range("J7").Value = 0.2979256989470644
$ws.Range("O7").Value = 0.003249135679578298
$ws.Range("P7").Value = 0.003249135679578299
$ws.Range("S7").Value = 0.0009680010183122096
$ws.Range("T7").Value = 0.0009680010183122097
# Row 8
$ws.Range("I8").Value = 0.2979256989470644
$ws.Range("J8").Value = 0.2979256989470644
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02449766666666667
$ws.Range("N8").Value = 0.073493
$ws.Range("O8").Value = 0.0007539688752111494
$ws.Range("P8").Value = 0.0007539688752111494
$ws.Range("Q8").Value = 0.6408449391687778
$ws.Range("R8").Value = 5.767604452519
$ws.Range("S8").Value = 0.0002246267041316137
$ws.Range("T8").Value = 0.0002246267041316137
# Row 9
$ws.Range("I9").Value = 0.2979256989470644
$ws.Range("J9").Value = 0.2979256989470644
$ws.Range("M9").Value = 32.353591
$ws.Range("N9").Value = 97.060773
$ws.Range("O9").Value = 0.9957520014958525
$ws.Range("P9").Value = 0.9957520014958525
$ws.Range("Q9").Value = 846.3514235214178
$ws.Range("R9").Value = 7617.16281169276
$ws.Range("S9").Value = 0.2966601110235902
$ws.Range("T9").Value = 0.2966601110235902
# Row 10
$ws.Range("G10").Value = 6.299630666666666
$ws.Range("H10").Value = 18.898892
$ws.Range("I10").Value = 0.07174552491706633
$ws.Range("J10").Value = 0.07174552491706633
$ws.Range("M10").Value = 0.007957
$ws.Range("N10").Value = 0.023871
$ws.Range("O10").Value = 0.0002448939493579708
$ws.Range("P10").Value = 0.0002448939493579708
$ws.Range("Q10").Value = 0.05012616121466667
$ws.Range("R10").Value = 0.451135450932
$ws.Range("S10").Value = 0.00001757004494570107
$ws.Range("T10").Value = 0.00001757004494570107
# Row 11
$ws.Range("G11").Value = 6.299630666666666
$ws.Range("H11").Value = 18.898892
$ws.Range("I11").Value = 0.07174552491706633
$ws.Range("J11").Value = 0.07174552491706633
$ws.Range("O11").Value = 0.003249135679578298
$ws.Range("P11").Value = 0.003249135679578299
$ws.Range("Q11").Value = 0.6650499096031112
$ws.Range("R11").Value = 5.985449186428
$ws.Range("S11").Value = 0.0002331109448581141
$ws.Range("T11").Value = 0.0002331109448581141
# Row 12
$ws.Range("G12").Value = 6.299630666666666
$ws.Range("H12").Value = 18.898892
$ws.Range("I12").Value = 0.07174552491706633
$ws.Range("J12").Value = 0.07174552491706633
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.02449766666666667
$ws.Range("N12").Value = 0.073493
$ws.Range("O12").Value = 0.0007539688752111494
$ws.Range("P12").Value = 0.0007539688752111494
$ws.Range("Q12").Value = 0.1543262521951111
$ws.Range("R12").Value = 1.388936269756
$ws.Range("S12").Value = 0.00005409389272315399
$ws.Range("T12").Value = 0.00005409389272315399
# Row 13
$ws.Range("G13").Value = 6.299630666666666
$ws.Range("H13").Value = 18.898892
$ws.Range("I13").Value = 0.07174552491706633
$ws.Range("J13").Value = 0.07174552491706633
$ws.Range("M13").Value = 32.353591
$ws.Range("N13").Value = 97.060773
$ws.Range("O13").Value = 0.9957520014958525
$ws.Range("P13").Value = 0.9957520014958525
$ws.Range("Q13").Value = 203.8156740403907
$ws.Range("R13").Value = 1834.341066363516
$ws.Range("S13").Value = 0.07144075003453935
$ws.Range("T13").Value = 0.07144075003453935
# Row 14
$ws.Range("G14").Value = 51.62686066666667
$ws.Range("H14").Value = 154.880582
$ws.Range("I14").Value = 0.5879703770491272
$ws.Range("J14").Value = 0.5879703770491272
$ws.Range("M14").Value = 0.007957
$ws.Range("N14").Value = 0.023871
$ws.Range("O14").Value = 0.0002448939493579708
$ws.Range("P14").Value = 0.0002448939493579708
$ws.Range("Q14").Value = 0.4107949303246667
$ws.Range("R14").Value = 3.697154372922
$ws.Range("S14").Value = 0.0001439903877410559
$ws.Range("T14").Value = 0.0001439903877410559
# Row 15
$ws.Range("G15").Value = 51.62686066666667
$ws.Range("H15").Value = 154.880582
$ws.Range("I15").Value = 0.5879703770491272
$ws.Range("J15").Value = 0.5879703770491272
$ws.Range("O15").Value = 0.003249135679578298
$ws.Range("P15").Value = 0.003249135679578299
$ws.Range("Q15").Value = 5.450230471626445
$ws.Range("R15").Value = 49.05207424463801
$ws.Range("S15").Value = 0.001910395530605424
$ws.Range("T15").Value = 0.001910395530605425
# Row 16
$ws.Range("G16").Value = 51.62686066666667
$ws.Range("H16").Value = 154.880582
$ws.Range("I16").Value = 0.5879703770491272
$ws.Range("J16").Value = 0.5879703770491272
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.02449766666666667
$ws.Range("N16").Value = 0.073493
$ws.Range("O16").Value = 0.0007539688752111494
$ws.Range("P16").Value = 0.0007539688752111494
$ws.Range("Q16").Value = 1.264737623658444
$ws.Range("R16").Value = 11.382638612926
$ws.Range("S16").Value = 0.0004433113638412059
$ws.Range("T16").Value = 0.0004433113638412059
# Row 17
$ws.Range("G17").Value = 51.62686066666667
$ws.Range("H17").Value = 154.880582
$ws.Range("I17").Value = 0.5879703770491272
$ws.Range("J17").Value = 0.5879703770491272
$ws.Range("M17").Value = 32.353591
$ws.Range("N17").Value = 97.060773
$ws.Range("O17").Value = 0.9957520014958525
$ws.Range("P17").Value = 0.9957520014958525
$ws.Range("Q17").Value = 1670.314334623321
$ws.Range("R17").Value = 15032.82901160989
$ws.Range("S17").Value = 0.5854726797669395
$ws.Range("T17").Value = 0.5854726797669395
